# Apply updated crypto price/volume figures (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" cells: force text so values like "1.00" / "0.998" keep their
# exact digits instead of Excel auto-converting the typed text to a number.
# ClearFormats() afterwards drops the temporary text format so the cell keeps
# the original (default) style, matching the source data.
function Set-PriceText($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-PriceText 'D2' '64.922.97'
$ws.Range('E2').Value = '  +2.29%  '
Set-PriceText 'D3' '2.639.64'
$ws.Range('E3').Value = '  +2.17%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-PriceText 'D5' '595.36'
$ws.Range('E5').Value = '  +0.99%  '
Set-PriceText 'D6' '155.34'
$ws.Range('E6').Value = '  +3.40%  '
Set-PriceText 'D7' '1.00'
$ws.Range('E7').Value = '  +0.04%  '
Set-PriceText 'D8' '0.591'
$ws.Range('E8').Value = '  +0.68%  '
Set-PriceText 'D9' '0.117'
$ws.Range('E9').Value = '  +5.58%  '
$ws.Range('E10').Value = '  +3.99%  '
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('E12').Value = '  +1.81%  '
Set-PriceText 'D13' '29.07'
$ws.Range('E13').Value = '  +5.26%  '
$ws.Range('E14').Value = '  +19.60%  '
Set-PriceText 'D15' '3.117.73'
$ws.Range('E15').Value = '  +2.27%  '
Set-PriceText 'D16' '64.844.54'
$ws.Range('E16').Value = '  +2.40%  '
Set-PriceText 'D17' '2.640.85'
$ws.Range('E17').Value = '  +2.04%  '
Set-PriceText 'D18' '12.55'
$ws.Range('E18').Value = '  +2.53%  '
$ws.Range('E19').Value = '  +1.57%  '
Set-PriceText 'D20' '351.44'
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('E21').Value = '  +6.12%  '
$ws.Range('E22').Value = '  +0.18%  '
Set-PriceText 'D23' '67.91'
$ws.Range('E23').Value = '  +0.93%  '
$ws.Range('E24').Value = '  -0.35%  '
Set-PriceText 'D25' '9.51'
$ws.Range('E25').Value = '  +3.97%  '
$ws.Range('E26').Value = '  -2.27%  '
Set-PriceText 'D27' '8.10'
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('E28').Value = '  +0.24%  '
Set-PriceText 'D29' '0.998'
$ws.Range('E29').Value = '  -0.16%  '
Set-PriceText 'D30' '0.0₃0936'
$ws.Range('E30').Value = '  +8.98%  '
Set-PriceText 'D31' '2.10'
$ws.Range('E31').Value = '  +2.75%  '
Set-PriceText 'D32' '511.13'
$ws.Range('E32').Value = '  -7.61%  '
Set-PriceText 'D33' '1.77'
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('E34').Value = '  +6.65%  '
Set-PriceText 'D35' '6.31'
$ws.Range('E35').Value = '  +4.09%  '
Set-PriceText 'D36' '0.425'
$ws.Range('E36').Value = '  +2.56%  '
Set-PriceText 'D37' '164.71'
$ws.Range('E37').Value = '  -1.11%  '
Set-PriceText 'D38' '20.19'
$ws.Range('E38').Value = '  +3.37%  '
$ws.Range('E39').Value = '  +4.37%  '
Set-PriceText 'D40' '1.00'
$ws.Range('E40').Value = '  +0.07%  '
Set-PriceText 'D41' '0.999'
$ws.Range('E41').Value = '  +0.02%  '
Set-PriceText 'D42' '42.22'
$ws.Range('E42').Value = '  +6.02%  '
Set-PriceText 'D43' '164.88'
$ws.Range('E43').Value = '  -0.86%  '
Set-PriceText 'D44' '4.09'
$ws.Range('E44').Value = '  +1.58%  '
$ws.Range('E45').Value = '  +3.30%  '
Set-PriceText 'D46' '22.89'
$ws.Range('E46').Value = '  -0.99%  '
$ws.Range('E47').Value = '  +3.57%  '
Set-PriceText 'D48' '0.647'
$ws.Range('E48').Value = '  +2.72%  '
Set-PriceText 'D49' '0.0254'
$ws.Range('E49').Value = '  +1.03%  '
Set-PriceText 'D50' '0.0980'
$ws.Range('E50').Value = '  +1.71%  '
Set-PriceText 'D51' '19.29'
$ws.Range('E51').Value = '  +0.62%  '
